# "passe sur la relecture de thierry de w6"
#
# 1) Slide 8, shape "ZoneTexte 5", 3rd paragraph: fix the typo
#    "bloque" -> "bloc" and split the run into three runs
#    ("On saute le " / "bloc " / "de code des classes englobantes"),
#    all keeping the original run formatting.
#
# 2) Slide 9, shape "Espace reserve du contenu 2" (#1), 9th paragraph:
#    merge the two runs "ins = C" + "()" back into a single run
#    "ins = C()", keeping the first run's formatting.

$p = $ppt.ActivePresentation

# --- Change 1 : slide 8 --------------------------------------------------
$s8  = $p.Slides.Item(8)
$sh8 = $s8.Shapes.Item(2)
$tr8 = $sh8.TextFrame.TextRange
$para3 = $tr8.Paragraphs(3, 1)

# Replace "bloque " (chars 13-19) with "bloc " so the run splits into
# three runs while each keeps the paragraph's original character formatting.
$sub = $para3.Characters(13, 7)
$sub.Text = "bloc "

# --- Change 2 : slide 9 --------------------------------------------------
$s9  = $p.Slides.Item(9)
$sh9 = $s9.Shapes.Item(1)
$tr9 = $sh9.TextFrame.TextRange
$para9 = $tr9.Paragraphs(9, 1)

$run1 = $para9.Runs(1)
$run1.Text = "ins = C()"

$run2 = $para9.Runs(2)
$run2.Text = ""
